# Turn the single "news blurb" slide into a 5-slide LinkedIn carousel draft:
#   Slide 1 - Hook            (existing slide, re-purposed)
#   Slide 2 - Key Highlights  (new)
#   Slide 3 - Why This Matters(new)
#   Slide 4 - Industry Impact (new)
#   Slide 5 - Source          (new)
#
# Every slide shares the same two-textbox recipe: a large "headline"
# textbox near the top and a smaller bulleted-paragraph textbox lower
# on the slide. Helper functions below build that recipe on demand.

$p = $ppt.ActivePresentation

function Add-HeadingBox($slide, [string]$text) {
    # 457200/274320/8229600/914400 EMU == 36/21.6/648/72 pt
    $tb = $slide.Shapes.AddTextbox(1, 36, 21.6, 648, 72)
    $tb.TextFrame.WordWrap = $false
    $tb.TextFrame.AutoSize = 1
    $tb.TextFrame.TextRange.Text = $text
    $tb.TextFrame.TextRange.Font.Size = 28
    $tb.Fill.Visible = $false
    # AutoSize shrink-to-fit already resized the box height for us;
    # put it back to the intended 72pt without touching Left/Top/Width
    # (re-assigning those drifts the EMU value by rounding).
    $tb.Height = 72
    return $tb
}

function Add-BodyBox($slide, [string[]]$lines) {
    # 731520/4114800/7772400/1828800 EMU == 57.6/324/612/144 pt
    $tb = $slide.Shapes.AddTextbox(1, 57.6, 324, 612, 144)
    $tb.TextFrame.WordWrap = $false
    $tb.TextFrame.AutoSize = 1
    $joined = "`n" + ($lines -join "`r")
    $tb.TextFrame.TextRange.Text = $joined
    $sub = $tb.TextFrame.TextRange.Characters(2, $tb.TextFrame.TextRange.Length - 1)
    $sub.Font.Size = 18
    $tb.Fill.Visible = $false
    $tb.Height = 144
    return $tb
}

# ---------------------------------------------------------------
# Slide 1 - reuse the existing slide, strip it down to the bare
# title placeholder, then add the new headline/body textboxes.
# ---------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# The "Content Placeholder" shape is a layout placeholder; the first
# Delete() only clears its content back to the empty inherited state,
# so it has to be called twice to actually remove the shape.
$s1.Shapes.Item("Content Placeholder 2").Delete()
$s1.Shapes.Item("Content Placeholder 2").Delete()

# The source-image picture is no longer part of the carousel look.
$s1.Shapes.Item("Picture 3").Delete()

# Clear the old headline text out of the title placeholder (kept as
# an empty placeholder, matching the other slides in the deck).
$s1.Shapes.Item("Title 1").TextFrame.TextRange.Text = ""

Add-HeadingBox $s1 "NTPC Green Energy board okays 50:50 JV with GAIL - India Infoline" | Out-Null
Add-BodyBox $s1 @("Strategic joint venture in India’s clean energy sector") | Out-Null

# ---------------------------------------------------------------
# Slides 2-5 - brand new "Title Only" slides, each with a headline
# textbox and a bulleted body textbox.
# ---------------------------------------------------------------
$s2 = $p.Slides.Add(2, 11)
Add-HeadingBox $s2 "Key Highlights" | Out-Null
Add-BodyBox $s2 @(
    "NTPC Green Energy approves a 50:50 JV with GAIL",
    "Focus on renewable and clean energy projects",
    "Strengthens public-sector collaboration"
) | Out-Null

$s3 = $p.Slides.Add(3, 11)
Add-HeadingBox $s3 "Why This Matters" | Out-Null
Add-BodyBox $s3 @(
    "Accelerates India’s energy transition",
    "Supports green hydrogen and renewables",
    "Enhances long-term energy security"
) | Out-Null

$s4 = $p.Slides.Add(4, 11)
Add-HeadingBox $s4 "Industry Impact" | Out-Null
Add-BodyBox $s4 @(
    "Boosts investor confidence in green energy",
    "Encourages large-scale clean infrastructure",
    "Aligns with India’s net-zero goals"
) | Out-Null

$s5 = $p.Slides.Add(5, 11)
Add-HeadingBox $s5 "Source" | Out-Null
Add-BodyBox $s5 @(
    "India Infoline",
    "Read full article online"
) | Out-Null

Write-Host "Slides:" $p.Slides.Count
